# Daily update at 8 AM UTC
#
# The sheet keeps a running daily tally (Day / Chase / Bryce / Zach).
# The previously-last row (row 55, 2024-12-14) was styled with the
# "latest row" date-only format; a new day's row is appended below it
# and that special formatting moves down to the new last row (row 56,
# 2024-12-15), while row 55 reverts to the regular date+time format
# used by every other data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the old "last row" (55) to the standard date/time number format.
$ws.Range("A55").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's row (56) with the next serial date and that day's values.
$ws.Range("A56").Value = 45641
$ws.Range("A56").NumberFormat = "YYYY-MM-DD"

$ws.Range("B56").Value = 135
$ws.Range("C56").Value = 120
$ws.Range("D56").Value = 128
